$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Append the new monthly data row (row 57: 2024-08-01 / serial 45505).
#    Copy row 56's formatting down first so the new date cell (A57) reuses
#    the existing "date" style instead of allocating a brand-new one, then
#    overwrite the copied values with the real data for the new row.
# ---------------------------------------------------------------------------
$ws.Cells.Item(56, 1).Copy($ws.Cells.Item(57, 1))

$ws.Cells.Item(57, 1).Value = 45505
$ws.Cells.Item(57, 2).Value = -0.439
$ws.Cells.Item(57, 3).Value = 0.422
$ws.Cells.Item(57, 4).Value = -0.979
$ws.Cells.Item(57, 5).Value = 0.35
$ws.Cells.Item(57, 6).Value = 1.737

# ---------------------------------------------------------------------------
# 2) Switch the Date column's number format from the custom "mm/dd/yyyy"
#    format to the built-in date format (numFmtId 14), for every date cell
#    in the column (A2:A57), including the row just added above.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 1).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(2, 1).Copy()
$ws.Range("A3:A57").PasteSpecial(-4122)
